# Insert a new "total_concentrations" worksheet right after
# "input_concentrations" (i.e. before "equilibrium_concentrations"),
# and populate it with the molecule1/molecule2 total-concentration data.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("input_concentrations")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "total_concentrations"

$newSheet.Range("A1").Value = "molecule1"
$newSheet.Range("B1").Value = "molecule2"

$newSheet.Range("A2").Value = 0.06098
$newSheet.Range("B2").Value = 0.0003999

$newSheet.Range("A3").Value = 0.06128
$newSheet.Range("B3").Value = 0.0006998

$newSheet.Range("A4").Value = 0.06098
$newSheet.Range("B4").Value = 0.000978

$newSheet.Range("A5").Value = 0.06208
$newSheet.Range("B5").Value = 0.0015

$newSheet.Range("A6").Value = 0.06199
$newSheet.Range("B6").Value = 0.00199

$newSheet.Range("A7").Value = 0.0009877
$newSheet.Range("B7").Value = 0.0009877

$newSheet.Range("A8").Value = 0.008865
$newSheet.Range("B8").Value = 0.008865
